# Daily IST report: add CSV/MD/XLSX
# Adds a new date column (2026-02-26) before the total_files/unique_days
# summary columns, shifting total_files -> M and unique_days -> N, and
# refreshes the per-row totals to include the new day's submissions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing "2026-02-25" (K) column's width as the template for
# the brand-new "2026-02-26" column, and insert a fresh column before the
# old L (total_files) column. Excel shifts total_files -> M and
# unique_days -> N automatically, carrying their header styles with them.
$dayColWidth = $ws.Columns("K:K").ColumnWidth

$ws.Columns("L:L").Insert()
$ws.Columns("L:L").ColumnWidth = $dayColWidth

# New header cell for the inserted date column (style already inherited
# as bold/centered from the Insert shift). Force text format first so the
# date-like string is stored as a literal label, not auto-converted to a
# date serial (matches how the other date-header cells read).
$ws.Range("L1").NumberFormat = "@"
$ws.Range("L1").Value = "2026-02-26"

# Per-row data: row, 2026-02-26 submissions, new total_files, new unique_days
$rows = @(
    @(2,1,5,5),
    @(3,1,4,4),
    @(4,1,5,5),
    @(5,1,6,6),
    @(6,1,4,4),
    @(7,1,6,6),
    @(8,1,6,6),
    @(9,1,6,6),
    @(10,1,5,5),
    @(11,1,5,5),
    @(12,1,6,6),
    @(13,1,5,5),
    @(14,1,6,6),
    @(15,1,4,4),
    @(16,1,4,4),
    @(17,1,6,6),
    @(18,1,6,6),
    @(19,1,4,4),
    @(20,1,5,5),
    @(21,0,0,0),
    @(22,0,0,0),
    @(23,1,1,1),
    @(24,1,5,5),
    @(25,0,4,4),
    @(26,1,3,3),
    @(27,1,7,5),
    @(28,0,0,0),
    @(29,0,0,0),
    @(30,0,0,0),
    @(31,1,6,6),
    @(32,1,6,6),
    @(33,1,6,6),
    @(34,1,6,6),
    @(35,1,5,5),
    @(36,1,5,5),
    @(37,1,3,3),
    @(38,0,0,0),
    @(39,0,0,0),
    @(40,1,6,6),
    @(41,1,6,6),
    @(42,1,6,6),
    @(43,0,21,1),
    @(44,1,44,3),
    @(45,1,3,3),
    @(46,1,6,6),
    @(47,1,6,6),
    @(48,0,0,0),
    @(49,1,6,6),
    @(50,1,5,5),
    @(51,1,1,1),
    @(52,0,2,2),
    @(53,1,5,5),
    @(54,0,0,0),
    @(55,1,6,6),
    @(56,0,0,0),
    @(57,1,5,5),
    @(58,1,5,5),
    @(59,0,3,3),
    @(60,1,6,6),
    @(61,0,1,1),
    @(62,0,0,0),
    @(63,1,3,3),
    @(64,1,6,6),
    @(65,0,0,0),
    @(66,0,0,0),
    @(67,0,2,2),
    @(68,0,0,0),
    @(69,0,0,0),
    @(70,1,6,6),
    @(71,0,1,1),
    @(72,0,0,0),
    @(73,1,3,3),
    @(74,1,17,5),
    @(75,1,1,1),
    @(76,1,4,4),
    @(77,0,0,0),
    @(78,0,0,0),
    @(79,1,3,3),
    @(80,1,5,5),
    @(81,0,4,4),
    @(82,1,6,6),
    @(83,1,3,3),
    @(84,0,1,1),
    @(85,0,0,0),
    @(86,1,3,3),
    @(87,0,0,0),
    @(88,0,2,2),
    @(89,0,0,0),
    @(90,1,2,2),
    @(91,0,0,0),
    @(92,0,0,0),
    @(93,1,3,3),
    @(94,1,4,4),
    @(95,1,33,4),
    @(96,0,1,1),
    @(97,0,0,0),
    @(98,1,2,2),
    @(99,0,0,0),
    @(100,0,2,2),
    @(101,1,6,6),
    @(102,0,0,0),
    @(103,1,6,6),
    @(104,0,12,1),
    @(105,1,6,6),
    @(106,1,6,6),
    @(107,0,0,0),
    @(108,1,5,5),
    @(109,0,0,0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 12).Value = $r[1]
    $ws.Cells.Item($rowNum, 13).Value = $r[2]
    $ws.Cells.Item($rowNum, 14).Value = $r[3]
}
